# "Change some code login error page and monthly report pages"
#
# For the monthly-report workbook, the author fixed a small typo in the
# "Lot No 1" column header (missing period after "No") and left the
# worksheet scrolled/selected over that column (T1) after making the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the "Lot No 1" header cell (column T, row 1) robustly instead of
# hard-coding the address, then correct the typo -> "Lot No. 1".
$target = $ws.Cells.Find("Lot No 1")
if ($target) {
    $target.Value = "Lot No. 1"
} else {
    $ws.Range("T1").Value = "Lot No. 1"
}

# Reflect the author's on-screen selection at the time of the edit.
$ws.Range("T1").Select()
